# Remove the "Title_Remark" column (column B) from the sample book list.
# This deletes the entire column, shifting Author/Publisher/Pub_Year/ISBN/
# Binding one column to the left (C->B, D->C, E->D, F->E, G->F) and
# shrinking the used range from A1:G24 to A1:F24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").EntireColumn.Delete()
